# Weekly update: insert 3 new price records at the top of the Espárragos
# block (pushing the existing history down by 3 rows) for the new
# reporting date 45225 (2023-10-26), keeping the newest-first ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 100; everything that was on rows
# 100-187 moves down to rows 103-190.
$ws.Range("A100:A102").EntireRow.Insert()

# --- New row 100: Banquete ---
$ws.Cells.Item(100, 1).Value2 = 9
$ws.Cells.Item(100, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(100, 3).Value2 = "Metropolitana"
$ws.Cells.Item(100, 4).Value2 = 45225
$ws.Cells.Item(100, 5).Value2 = 13
$ws.Cells.Item(100, 6).Value2 = 300000000
$ws.Cells.Item(100, 7).Value2 = "Espárragos"
$ws.Cells.Item(100, 8).Value2 = "Sin especificar"
$ws.Cells.Item(100, 9).Value2 = "Banquete"
$ws.Cells.Item(100, 10).Value2 = 70
$ws.Cells.Item(100, 11).Value2 = 16000
$ws.Cells.Item(100, 12).Value2 = 16000
$ws.Cells.Item(100, 13).Value2 = 16000
$ws.Cells.Item(100, 14).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(100, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(100, 16).Value2 = 1600
$ws.Cells.Item(100, 17).Value2 = 10
$ws.Cells.Item(100, 18).Value2 = "Hortaliza"

# --- New row 101: Primera ---
$ws.Cells.Item(101, 1).Value2 = 9
$ws.Cells.Item(101, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(101, 3).Value2 = "Metropolitana"
$ws.Cells.Item(101, 4).Value2 = 45225
$ws.Cells.Item(101, 5).Value2 = 13
$ws.Cells.Item(101, 6).Value2 = 300000000
$ws.Cells.Item(101, 7).Value2 = "Espárragos"
$ws.Cells.Item(101, 8).Value2 = "Sin especificar"
$ws.Cells.Item(101, 9).Value2 = "Primera"
$ws.Cells.Item(101, 10).Value2 = 160
$ws.Cells.Item(101, 11).Value2 = 14000
$ws.Cells.Item(101, 12).Value2 = 14000
$ws.Cells.Item(101, 13).Value2 = 14000
$ws.Cells.Item(101, 14).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(101, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(101, 16).Value2 = 1400
$ws.Cells.Item(101, 17).Value2 = 10
$ws.Cells.Item(101, 18).Value2 = "Hortaliza"

# --- New row 102: Segunda ---
$ws.Cells.Item(102, 1).Value2 = 9
$ws.Cells.Item(102, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(102, 3).Value2 = "Metropolitana"
$ws.Cells.Item(102, 4).Value2 = 45225
$ws.Cells.Item(102, 5).Value2 = 13
$ws.Cells.Item(102, 6).Value2 = 300000000
$ws.Cells.Item(102, 7).Value2 = "Espárragos"
$ws.Cells.Item(102, 8).Value2 = "Sin especificar"
$ws.Cells.Item(102, 9).Value2 = "Segunda"
$ws.Cells.Item(102, 10).Value2 = 97
$ws.Cells.Item(102, 11).Value2 = 12000
$ws.Cells.Item(102, 12).Value2 = 12000
$ws.Cells.Item(102, 13).Value2 = 12000
$ws.Cells.Item(102, 14).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(102, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(102, 16).Value2 = 1200
$ws.Cells.Item(102, 17).Value2 = 10
$ws.Cells.Item(102, 18).Value2 = "Hortaliza"
